# Existing instructor olayi handle edildi
# Split instructors that teach courses from more than one department
# (or more than one course) into their own separate rows on the
# "Instructors" sheet, and rename the COURSES column header to TEACHING.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Instructors")

# Rename header E1: COURSES -> TEACHING
$ws.Range("E1").Value = "TEACHING"

# --- Serdar Tasiran (COMP) -----------------------------------------
# Row 2 keeps Serdar Tasiran but now only references COMP 302 / Java.
$ws.Range("E2").Value = "COMP 302"
$ws.Range("F2").Value = "Java"

# Row 3 becomes a second row for Serdar Tasiran (COMP 131 / GRASP),
# copying formatting from row 2 first so the new cells inherit the
# right styles/borders.
$ws.Range("A2:F2").Copy()
$ws.Range("A3:F3").PasteSpecial(-4122)
$ws.Range("G2:J2").Copy()
$ws.Range("G3:J3").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws.Range("A3").Value = "Serdar"
$ws.Range("B3").Value = "Tasiran"
$ws.Range("C3").Value = "mpinar@ku.edu.tr"
$ws.Range("D3").Value = "COMP"
$ws.Range("E3").Value = "COMP 131"
$ws.Range("F3").Value = "GRASP"
$ws.Rows.Item(3).RowHeight = 20

# --- Ozgur Baris Akan (ELEC) ----------------------------------------
# Rows 4 and 5 are currently blank placeholder rows; they become the
# two split rows for Ozgur Baris Akan (ELEC 201 / Matlab and
# ELEC 204 / Matlab).
$ws.Range("A2:F2").Copy()
$ws.Range("A4:F4").PasteSpecial(-4122)
$ws.Range("A4:F4").Copy()
$ws.Range("A5:F5").PasteSpecial(-4122)
$ws.Range("G2:J2").Copy()
$ws.Range("G4:J4").PasteSpecial(-4122)
$ws.Range("G4:J4").Copy()
$ws.Range("G5:J5").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws.Range("A4").Value = "Ozgur Baris"
$ws.Range("B4").Value = "Akan"
$ws.Range("C4").Value = "doztreves@ku.edu.tr"
$ws.Range("D4").Value = "ELEC"
$ws.Range("E4").Value = "ELEC 201"
$ws.Range("F4").Value = "Matlab"
$ws.Rows.Item(4).RowHeight = 20

$ws.Range("A5").Value = "Ozgur Baris"
$ws.Range("B5").Value = "Akan"
$ws.Range("C5").Value = "doztreves@ku.edu.tr"
$ws.Range("D5").Value = "ELEC"
$ws.Range("E5").Value = "ELEC 204"
$ws.Range("F5").Value = "Matlab"
$ws.Rows.Item(5).RowHeight = 15.65

# --- Extra trailing blank row -----------------------------------------
# A new blank row 11 is appended (matching the dimension growing from
# A1:J10 to A1:J11); copy formatting from the existing blank row 10.
$ws.Range("A10:J10").Copy()
$ws.Range("A11:J11").PasteSpecial(-4122)
$excel.CutCopyMode = 0
$ws.Rows.Item(11).RowHeight = 15.65
$ws.Range("A11:J11").ClearContents()

# --- Hyperlinks on the MAIL column -----------------------------------
# The engine only supports a full rebuild of the hyperlink collection
# (per-cell delete removes every link on the sheet), so wipe everything
# and re-add the four mailto links in the right order: C2/C3 -> mpinar,
# C4/C5 -> doztreves.
$ws.Range("C2").Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("C2"), "mailto:mpinar@ku.edu.tr")
$ws.Hyperlinks.Add($ws.Range("C3"), "mailto:mpinar@ku.edu.tr")
$ws.Hyperlinks.Add($ws.Range("C4"), "mailto:doztreves@ku.edu.tr")
$ws.Hyperlinks.Add($ws.Range("C5"), "mailto:doztreves@ku.edu.tr")

# Adding a hyperlink auto-applies Excel's built-in "Hyperlink" cell
# style (underline + colored font); the source file keeps the plain
# data-row formatting instead, so copy it back from column D (plain,
# unstyled) onto the four mail cells without touching their values.
$ws.Range("D2").Copy()
$ws.Range("C2").PasteSpecial(-4122)
$ws.Range("D3").Copy()
$ws.Range("C3").PasteSpecial(-4122)
$ws.Range("D4").Copy()
$ws.Range("C4").PasteSpecial(-4122)
$ws.Range("D5").Copy()
$ws.Range("C5").PasteSpecial(-4122)
$excel.CutCopyMode = 0
